$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 7001.5
$ws.Range("I76").Value = 7001.5
$ws.Range("K76").Value = 7001.5
$ws.Range("M76").Value = -6686.5

$ws.Range("H79").Value = 7001.5
$ws.Range("I79").Value = 7001.5
$ws.Range("K79").Value = 7001.5
$ws.Range("M79").Value = -5909.5

$ws.Range("H80").Value = 865.6070999999999
$ws.Range("I80").Value = 695.7857
$ws.Range("J80").Value = 1035.4286
$ws.Range("K80").Value = 2087.3571
$ws.Range("L80").Value = 3106.2858
$ws.Range("M80").Value = -1089.3571
$ws.Range("N80").Value = -5102.2858

$ws.Range("H83").Value = 865.6070999999999
$ws.Range("I83").Value = 695.7857
$ws.Range("J83").Value = 1035.4286
$ws.Range("K83").Value = 6262.071300000001
$ws.Range("L83").Value = 9318.857399999999
$ws.Range("M83").Value = -1270.071300000001
$ws.Range("N83").Value = -19302.8574

$ws.Range("H94").Value = 13958.053
$ws.Range("I94").Value = 14177.944
$ws.Range("K94").Value = 14177.944
$ws.Range("M94").Value = -13726.944

$ws.Range("H135").Value = 1116.5714
$ws.Range("I135").Value = 1116.5714
$ws.Range("K135").Value = 10049.1426
$ws.Range("M135").Value = -7514.142600000001

$ws.Range("H137").Value = 6630.269
$ws.Range("J137").Value = 10844.286
$ws.Range("L137").Value = 32532.858
$ws.Range("N137").Value = -37632.858

$ws.Range("H138").Value = 3073.7144
$ws.Range("I138").Value = 1347.069
$ws.Range("J138").Value = 3881.3386
$ws.Range("K138").Value = 4041.207
$ws.Range("L138").Value = 11644.0158
$ws.Range("M138").Value = 1098.793
$ws.Range("N138").Value = -21924.0158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 6499.857
$ws.Range("I25").Value = 1799.5
$ws.Range("J25").Value = 8380
$ws.Range("K25").Value = 1799.5
$ws.Range("L25").Value = 8380
$ws.Range("M25").Value = -1397.5
$ws.Range("N25").Value = -9184

$ws.Range("H32").Value = 2711.8
$ws.Range("I32").Value = 2238.279
$ws.Range("J32").Value = 5620.5713
$ws.Range("K32").Value = 2238.279
$ws.Range("L32").Value = 5620.5713
$ws.Range("M32").Value = -1951.279
$ws.Range("N32").Value = -6194.5713

$ws.Range("H45").Value = 24222.637
$ws.Range("I45").Value = 18103.889
$ws.Range("K45").Value = 18103.889
$ws.Range("M45").Value = -17726.889

$ws.Range("H74").Value = 14144.223
$ws.Range("J74").Value = 35416.332
$ws.Range("L74").Value = 35416.332
$ws.Range("N74").Value = -37164.332

$ws.Range("H77").Value = 14144.223
$ws.Range("J77").Value = 35416.332
$ws.Range("L77").Value = 177081.66
$ws.Range("N77").Value = -185817.66

$ws.Range("H110").Value = 7189.607
$ws.Range("I110").Value = 8814.5
$ws.Range("J110").Value = 5564.7144
$ws.Range("K110").Value = 8814.5
$ws.Range("L110").Value = 5564.7144
$ws.Range("M110").Value = -6769.5
$ws.Range("N110").Value = -9654.714400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3441.543
$ws.Range("I20").Value = 3156.4827
$ws.Range("J20").Value = 4819.3335
$ws.Range("K20").Value = 3156.4827
$ws.Range("L20").Value = 4819.3335
$ws.Range("M20").Value = -2909.4827
$ws.Range("N20").Value = -5313.3335

$ws.Range("H41").Value = 600684
$ws.Range("J41").Value = 600684
$ws.Range("L41").Value = 600684
$ws.Range("N41").Value = -601460

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5406137.5
$ws.Range("I31").Value = 6970826
$ws.Range("J31").Value = 41490.855
$ws.Range("K31").Value = 6970826
$ws.Range("L31").Value = 41490.855
$ws.Range("M31").Value = -6970531
$ws.Range("N31").Value = -42080.855

$ws.Range("H34").Value = 5406137.5
$ws.Range("I34").Value = 6970826
$ws.Range("J34").Value = 41490.855
$ws.Range("K34").Value = 6970826
$ws.Range("L34").Value = 41490.855
$ws.Range("M34").Value = -6970624
$ws.Range("N34").Value = -41894.855

$ws.Range("H99").Value = 3916.5334
$ws.Range("I99").Value = 3756.5
$ws.Range("J99").Value = 4556.6665
$ws.Range("K99").Value = 3756.5
$ws.Range("L99").Value = 4556.6665
$ws.Range("M99").Value = -2258.5
$ws.Range("N99").Value = -7552.6665

$ws.Range("H126").Value = 3916.5334
$ws.Range("I126").Value = 3756.5
$ws.Range("J126").Value = 4556.6665
$ws.Range("K126").Value = 11269.5
$ws.Range("L126").Value = 13669.9995
$ws.Range("M126").Value = -8799.5
$ws.Range("N126").Value = -18609.9995

$ws.Range("H132").Value = 5001.2666
$ws.Range("I132").Value = 4820.727
$ws.Range("J132").Value = 5497.75
$ws.Range("K132").Value = 14462.181
$ws.Range("L132").Value = 16493.25
$ws.Range("M132").Value = -11932.181
$ws.Range("N132").Value = -21553.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 9390
$ws.Range("I87").Value = 6813.75
$ws.Range("K87").Value = 20441.25
$ws.Range("M87").Value = -19193.25

$ws.Range("H90").Value = 9390
$ws.Range("I90").Value = 6813.75
$ws.Range("K90").Value = 61323.75
$ws.Range("M90").Value = -55083.75

$ws.Range("H113").Value = 533.2
$ws.Range("I113").Value = 528.4286
$ws.Range("J113").Value = 600
$ws.Range("K113").Value = 1585.2858
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = 584.7142000000001
$ws.Range("N113").Value = -6140

$ws.Range("H140").Value = 3137.1
$ws.Range("I140").Value = 2896.5881
$ws.Range("K140").Value = 8689.764299999999
$ws.Range("M140").Value = -3509.764299999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2002.7858
$ws.Range("I102").Value = 2002.7858
$ws.Range("K102").Value = 2002.7858
$ws.Range("M102").Value = -380.7858000000001

$ws.Range("H122").Value = 590
$ws.Range("J122").Value = 587.5
$ws.Range("L122").Value = 1762.5
$ws.Range("N122").Value = -6662.5

$ws.Range("H126").Value = 21481.646
$ws.Range("I126").Value = 28265.666
$ws.Range("K126").Value = 84796.99800000001
$ws.Range("M126").Value = -82326.99800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8971.65
$ws.Range("I7").Value = 9385.166999999999
$ws.Range("K7").Value = 9385.166999999999
$ws.Range("M7").Value = -9273.166999999999

$ws.Range("H61").Value = 5192.357
$ws.Range("J61").Value = 6166.3335
$ws.Range("L61").Value = 6166.3335
$ws.Range("N61").Value = -6570.3335

$ws.Range("H113").Value = 5192.357
$ws.Range("J113").Value = 6166.3335
$ws.Range("L113").Value = 6166.3335
$ws.Range("N113").Value = -10506.3335

$ws.Range("H126").Value = 8971.65
$ws.Range("I126").Value = 9385.166999999999
$ws.Range("K126").Value = 28155.501
$ws.Range("M126").Value = -25685.501

$ws.Range("H132").Value = 6645.4546
$ws.Range("J132").Value = 7500
$ws.Range("L132").Value = 22500
$ws.Range("N132").Value = -27560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 22000
$ws.Range("J39").Value = 22000
$ws.Range("L39").Value = 22000
$ws.Range("N39").Value = -22826

$ws.Range("H107").Value = 1159.4584
$ws.Range("I107").Value = 956.9
$ws.Range("J107").Value = 1304.1428
$ws.Range("K107").Value = 2870.7
$ws.Range("L107").Value = 3912.4284
$ws.Range("M107").Value = -950.6999999999998
$ws.Range("N107").Value = -7752.428400000001

$ws.Range("H126").Value = 3427.4443
$ws.Range("I126").Value = 3128.35
$ws.Range("K126").Value = 9385.049999999999
$ws.Range("M126").Value = -6915.049999999999

$ws.Range("H135").Value = 52885
$ws.Range("J135").Value = 52885
$ws.Range("L135").Value = 52885
$ws.Range("N135").Value = -63025

$ws.Range("H136").Value = 2988.4614
$ws.Range("I136").Value = 2715
$ws.Range("J136").Value = 3900
$ws.Range("K136").Value = 8145
$ws.Range("L136").Value = 11700
$ws.Range("M136").Value = -5595
$ws.Range("N136").Value = -16800
